$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '42.654.74'
$ws.Range('E2').Value = '  +0.77%  '
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '2.519.10'
$ws.Range('E3').Value = '  +0.45%  '
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '318.04'
$ws.Range('E5').Value = '  +4.16%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '95.62'
$ws.Range('E6').Value = '  -0.24%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.575'
$ws.Range('E7').Value = '  -1.35%  '
$ws.Range('E8').Value = '  -0.13%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.532'
$ws.Range('E9').Value = '  -0.75%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '35.77'
$ws.Range('E10').Value = '  -1.78%  '
$ws.Range('E11').Value = '  +0.24%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '7.55'
$ws.Range('E12').Value = '  -1.27%  '
$ws.Range('E13').Value = '  -3.36%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '2.905.41'
$ws.Range('E14').Value = '  +0.49%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '2.528.15'
$ws.Range('E15').Value = '  +1.88%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '15.31'
$ws.Range('E16').Value = '  -0.60%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '0.851'
$ws.Range('E17').Value = '  -0.52%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '42.710.89'
$ws.Range('E18').Value = '  +0.89%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '12.85'
$ws.Range('E19').Value = '  -0.58%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '6.71'
$ws.Range('E20').Value = '  +4.55%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '0.0₃0959'
$ws.Range('E21').Value = '  -1.31%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '69.67'
$ws.Range('E22').Value = '  -2.21%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '249.77'
$ws.Range('E23').Value = '  -1.19%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '2.96'
$ws.Range('E24').Value = '  +1.73%  '
$ws.Range('E25').Value = '  +3.49%  '
$ws.Range('B26').Value = 'EthereumClassic'
$ws.Range('C26').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '26.57'
$ws.Range('E26').Value = '  -1.47%  '
$ws.Range('B27').Value = 'Dai'
$ws.Range('C27').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '1.00'
$ws.Range('E27').Value = '  +0.00%  '
$ws.Range('B28').Value = 'Toncoin'
$ws.Range('C28').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '2.43'
$ws.Range('E28').Value = '  +4.41%  '
$ws.Range('B29').Value = 'InjectiveProtocol'
$ws.Range('C29').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '41.57'
$ws.Range('E29').Value = '  +11.26%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '10.26'
$ws.Range('E30').Value = '  +1.12%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '5.96'
$ws.Range('E31').Value = '  +0.68%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '157.64'
$ws.Range('E32').Value = '  +1.89%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '2.14'
$ws.Range('E33').Value = '  +3.33%  '
$ws.Range('B34').Value = 'Celestia'
$ws.Range('C34').Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '19.27'
$ws.Range('E34').Value = '  +0.48%  '
$ws.Range('B35').Value = 'WEMIXToken'
$ws.Range('C35').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '2.67'
$ws.Range('E35').Value = '  +2.73%  '
$ws.Range('E36').Value = '  -0.19%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '0.0778'
$ws.Range('E37').Value = '  -0.82%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.111'
$ws.Range('E38').Value = '  -2.49%  '
$ws.Range('E39').Value = '  -0.72%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '23.45'
$ws.Range('E40').Value = '  -2.95%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '2.30'
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '0.0305'
$ws.Range('E42').Value = '  +1.72%  '
$ws.Range('E43').Value = '  +0.36%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '3.78'
$ws.Range('E44').Value = '  -1.16%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '3.32'
$ws.Range('E45').Value = '  -1.88%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '2.030.78'
$ws.Range('E46').Value = '  +0.02%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '84.51'
$ws.Range('E47').Value = '  +0.20%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '8.92'
$ws.Range('E48').Value = '  -0.86%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '75.33'
$ws.Range('E49').Value = '  +3.26%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '2.756.88'
$ws.Range('E50').Value = '  +0.27%  '
$ws.Range('E51').Value = '  +2.54%  '
